# Weekly update: insert two new observation rows (Cilantro, Vega Modelo de
# Temuco) at the top of the data block (row 376) and push the existing
# 85 rows (old 376..460) down by two rows so they end up at 378..462.
# Excel's native row-insert semantics take care of re-numbering every row
# reference below the insertion point (and bump the used-range dimension
# from R460 to R462) exactly the way a human editor pasting two new rows
# at the top of the table would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 376 - this shifts rows
# 376:460 down to 378:462, matching the rest of the diff automatically.
$ws.Rows("376:377").Insert()

# Row 376 - brand new observation
$ws.Cells.Item(376, 1).Value  = 10
$ws.Cells.Item(376, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(376, 3).Value  = "La Araucanía"
$ws.Cells.Item(376, 4).Value  = 44855
$ws.Cells.Item(376, 5).Value  = 9
$ws.Cells.Item(376, 6).Value  = 100112040
$ws.Cells.Item(376, 7).Value  = "Cilantro"
$ws.Cells.Item(376, 8).Value  = "Sin especificar"
$ws.Cells.Item(376, 9).Value  = "Primera"
$ws.Cells.Item(376, 10).Value = 40
$ws.Cells.Item(376, 11).Value = 6500
$ws.Cells.Item(376, 12).Value = 6500
$ws.Cells.Item(376, 13).Value = 6500
$ws.Cells.Item(376, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(376, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(376, 16).Value = 3250
$ws.Cells.Item(376, 17).Value = 2
$ws.Cells.Item(376, 18).Value = "Hortaliza"

# Row 377 - brand new observation
$ws.Cells.Item(377, 1).Value  = 10
$ws.Cells.Item(377, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(377, 3).Value  = "La Araucanía"
$ws.Cells.Item(377, 4).Value  = 44855
$ws.Cells.Item(377, 5).Value  = 9
$ws.Cells.Item(377, 6).Value  = 100112040
$ws.Cells.Item(377, 7).Value  = "Cilantro"
$ws.Cells.Item(377, 8).Value  = "Sin especificar"
$ws.Cells.Item(377, 9).Value  = "Primera"
$ws.Cells.Item(377, 10).Value = 40
$ws.Cells.Item(377, 11).Value = 3000
$ws.Cells.Item(377, 12).Value = 3000
$ws.Cells.Item(377, 13).Value = 3000
$ws.Cells.Item(377, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(377, 15).Value = "Región Metropolitana"
$ws.Cells.Item(377, 16).Value = 1500
$ws.Cells.Item(377, 17).Value = 2
$ws.Cells.Item(377, 18).Value = "Hortaliza"
